# 2021 Q2 Quarterly Report WIP - Second draft for James
# Update Table 5 figures to reflect revised counts/percentages.

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "37 (30.6)" "37 (29.8)"
Replace-Text "27 (22.3)" "28 (22.6)"
Replace-Text "24 (19.8)" "26 (21.0)"
Replace-Text "12 (9.9)" "12 (9.7)"
Replace-Text "9 (7.4)" "9 (7.3)"
Replace-Text "2 (1.7)" "2 (1.6)"
